$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts the existing
# B column (the "a" picture / endImage) to C, and the existing
# C column (the "AB"/"a"/"b" label) to D, matching the diff.
$ws.Columns("B:B").Insert()

# Populate the new column B ("midImage") with the mid-point image
# path for each pair of rows (copy of the "a" image), except for
# row 25 which gets the new discrete opacity step image.
$ws.Range("B1").Value = "midImage"
$ws.Range("B2").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\1a.png"
$ws.Range("B3").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\1a.png"
$ws.Range("B4").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\2a.png"
$ws.Range("B5").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\2a.png"
$ws.Range("B6").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\3a.png"
$ws.Range("B7").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\3a.png"
$ws.Range("B8").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\4a.png"
$ws.Range("B9").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\4a.png"
$ws.Range("B10").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\5a.png"
$ws.Range("B11").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\5a.png"
$ws.Range("B12").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\6a.png"
$ws.Range("B13").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\6a.png"
$ws.Range("B14").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\7a.png"
$ws.Range("B15").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\7a.png"
$ws.Range("B16").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\8a.png"
$ws.Range("B17").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\8a.png"
$ws.Range("B18").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\9a.png"
$ws.Range("B19").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\9a.png"
$ws.Range("B20").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\10a.png"
$ws.Range("B21").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\10a.png"
$ws.Range("B22").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\11a.png"
$ws.Range("B23").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\11a.png"
$ws.Range("B24").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\12a.png"
$ws.Range("B25").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\12bapng"
$ws.Range("B26").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\13a.png"
$ws.Range("B27").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\13a.png"
$ws.Range("B28").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\14a.png"
$ws.Range("B29").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\14a.png"
$ws.Range("B30").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\15a.png"
$ws.Range("B31").Value = "C:\Users\hesselmann\Desktop\PsychoPyUniformityIllusion\UniformityIllusionPsychoPy\Pics\15a.png"

# Adjust column widths: B and C are now both 93 (character-width) wide.
# The Excel COM ColumnWidth property and the stored OOXML column width
# differ by a small constant offset (~5/6 of a character) in this
# runtime's unit conversion, so back that offset out to land on an
# exact stored width of 93.
$ws.Columns("B:B").ColumnWidth = 92.16666666666667
$ws.Columns("C:C").ColumnWidth = 92.16666666666667

# Update the view: scroll back to the top and select A10.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A10").Select()
